$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: column A (index), B, C, D for rows 2..38 (sheet row = index + 2)
$data = @(
    @(0, 1.572542121719036, 3.573675998788179, 0.1682440689762209),
    @(1, 2.471875601886715, 13.66778973488689, 0.2481175178011985),
    @(2, 6.869509479688297, 9.668529676418325, 0.2485055396418034),
    @(3, 11.86397520548146, 2.282088441392292, 0.2781697875166194),
    @(4, 14.46411246506577, 7.213386370790182, 0.2058729862662284),
    @(5, 21.39062923347744, 8.780482643419424, 0.2619635363017802),
    @(6, 21.89360964428272, 5.003452587337643, 0.3359027771647253),
    @(7, 22.27162400124058, 6.205759828660251, 0.56317255527414),
    @(8, 23.41847818429334, 2.462080416484365, 0.1905477177590064),
    @(9, 24.36625815205676, 8.07945614204414, 0.1397689782803794),
    @(10, 26.76731041766614, 8.442973860406552, 0.3434554275025283),
    @(11, 28.37458172113914, 3.190517626457088, 0.2946709233228973),
    @(12, 29.8140107588996, 7.929744089259489, 0.4838964756855428),
    @(13, 31.46021936588081, 2.111358546182495, 0.1753943948996613),
    @(14, 32.09965772108487, 7.711276599561025, 0.3068185749288009),
    @(15, 36.52588126029993, 3.26788865765766, 0.2063874415691686),
    @(16, 39.03568479388384, 6.674947261798255, 0.1613758054733849),
    @(17, 42.25720512618021, 6.941587135055538, 0.3783521117874042),
    @(18, 42.57807981599107, 2.292900466152926, 0.06930934494196311),
    @(19, 46.94562046946817, 11.80053932980652, 0.3892941784961114),
    @(20, 47.21446559881394, 5.569532801176738, 0.1966997888237451),
    @(21, 51.27594126760494, 3.013274111135734, 0.2460867459971646),
    @(22, 56.92549195466314, 3.395151305512115, 0.4654007084820373),
    @(23, 57.54322546582416, 7.974175543466883, 0.5925928358199157),
    @(24, 59.76555762766493, 3.459933374335276, 0.1801889967981982),
    @(25, 62.6385129131152, 6.766661366072175, 0.3554766858047512),
    @(26, 72.34714087577925, 6.576065042343197, 0.2873903724427079),
    @(27, 72.81327866107455, 8.947935788514114, 0.2090275889836607),
    @(28, 75.09401473901897, 5.539062107282198, 0.4191242911199765),
    @(29, 78.03653577836182, 7.091164606586593, 0.2187821653859486),
    @(30, 78.38626867725668, 8.546091232639348, 0.3841338768466209),
    @(31, 79.55872701124409, 7.760057732503734, 0.3799961378308657),
    @(32, 94.15225290176268, 8.027181716944559, 0.3165531640635245),
    @(33, 94.30932985152707, 12.33461981220447, 0.203206989039628),
    @(34, 94.75017836606995, 8.404860813064992, 0.1477723881304603),
    @(35, 98.37866363835198, 9.703326958923068, 0.4888594267858188),
    @(36, 98.96595071800607, 5.307800015201298, 0.1462619861178787)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $entry = $data[$i]

    if ($row -gt 33) {
        # New rows: clone the formatting used by the existing A-column cells
        # (bold, centered, thin-bordered) before writing the value.
        $ws.Range("A33").Copy()
        $ws.Range("A$row").PasteSpecial(-4122)  # xlPasteFormats
    }

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
